# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# This script adds a new "2509" period for both workers (mirroring the
# existing alternating CC/1047416130/YESENIA - CC/1128063438/OSMIRO pattern),
# updates the totals (Valor Mora / Cant. Periodos) accordingly, and moves the
# signature block down to make room for the two new data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the trailing signature rows (28,29) down to (30,31) to make room
#     for the two new worker/period rows being inserted after row 23. -----
$ws.Rows.Item(29).Copy()
$ws.Rows.Item(31).PasteSpecial(-4104)
$ws.Rows.Item(28).Copy()
$ws.Rows.Item(30).PasteSpecial(-4104)

$ws.Rows.Item(28).ClearContents()
$ws.Rows.Item(29).ClearContents()

# --- Insert the two new rows for period 2509, copying formatting/values
#     from the prior period's rows (22 -> 24, 23 -> 25). ------------------
$ws.Rows.Item(22).Copy()
$ws.Rows.Item(24).PasteSpecial(-4104)
$ws.Rows.Item(23).Copy()
$ws.Rows.Item(25).PasteSpecial(-4104)

$ws.Range("E24").Value = "2509"
$ws.Range("E25").Value = "2509"

# --- Update the summary figures: one more period, one more period's worth
#     of mora value (2 trabajadores x 64000). ------------------------------
$ws.Range("F13").Value = 5
$ws.Range("E11").Value = 640000

$wb.Save()
